# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
# Mirrors the original author's fix: the previous scraper only pulled team
# statistics, not the season W/L/T record, so three new columns are appended
# right after the existing last column (AC) with the season record repeated
# for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1, which
# carries the bold/centered/bordered header formatting) onto the three new
# header cells so they look consistent with the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# 2009 Houston Astros season record: 74 wins, 88 losses, 0 ties.
$wins = 74
$losses = 88
$ties = 0

$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
